$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Update the username/email in the test-data row to the latest challenge
# address (Webtest run), mirroring the prior "all assertions completed"
# edits to this sheet.
$ws.Range("B2").Value = "hf_challenge_1578862413072@hf413072.com"
